$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.456.63'
$ws.Range("E2").Value = '  +9.02%  '

$ws.Range("D3").Value = '1.679.57'
$ws.Range("E3").Value = '  +4.84%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.68'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9982'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.60%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3716'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.77%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3448'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.33%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '48.17'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +12.61%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.185'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.74%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07278'
$ws.Range("D11").Style = "Normal"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.000'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.14%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.44'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.87%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.139'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.28%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.753'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.72%  '

$ws.Range("D16").Value = '1.674.89'
$ws.Range("E16").Value = '  +4.62%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001110'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.16%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.9979'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.60%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06722'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.90%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '81.37'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.26%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.48'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.20%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.111'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.28%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.99'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.35%  '

$ws.Range("D24").Value = '24.394.41'
$ws.Range("E24").Value = '  +8.71%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.439'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.61%  '

$ws.Range("B26").Value = 'LidoDAOToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.674'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.55%  '

$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '152.62'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.49%  '

$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.57'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.32%  '

$ws.Range("B29").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C29").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D29").Value = '1.860.83'
$ws.Range("E29").Value = '  +4.51%  '

$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.83'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.18%  '

$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.357'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.26%  '

$ws.Range("B32").Value = 'HuobiToken'
$ws.Range("C32").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.033'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.99%  '

$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.9750'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.34%  '

$ws.Range("B34").Value = 'Stellar'
$ws.Range("C34").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08473'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.31%  '

$ws.Range("B35").Value = 'WEMIXTOKEN'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.672'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.03%  '

$ws.Range("B36").Value = 'Aptos'
$ws.Range("C36").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.45'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +4.78%  '

$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06508'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.33%  '

$ws.Range("B38").Value = 'FraxShare'
$ws.Range("C38").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.952'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.15%  '

$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.358'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.10%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02344'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.47%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.268'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.89%  '

$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2116'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.23%  '

$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6202'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.77%  '

$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9979'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.66%  '

$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.775'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.61%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5959'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.52%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.98'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.88%  '

$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '127.20'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.26%  '

$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.036'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.32%  '

$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07223'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.95%  '

$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '75.76'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.53%  '
